# Update cryptocurrency price/volume data per the Mon Jan 1 17:46:26 UTC 2024
# GitHub Actions refresh. Column map: B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }

$data = @{
    2 = @{ D="42.735.33"; E="  +0.44%  " }
    3 = @{ D="2.311.06"; E="  +0.69%  " }
    4 = @{ D="0.998"; E="  -0.35%  " }
    5 = @{ D="311.17"; E="  -1.60%  " }
    6 = @{ D="106.73"; E="  +2.34%  " }
    7 = @{ D="0.621"; E="  -1.03%  " }
    8 = @{ E="  -0.06%  " }
    9 = @{ E="  +0.60%  " }
    10 = @{ D="40.02"; E="  +0.92%  " }
    11 = @{ D="0.0914"; E="  +0.72%  " }
    12 = @{ D="8.38"; E="  -1.53%  " }
    13 = @{ E="  -1.60%  " }
    14 = @{ D="0.988"; E="  -1.66%  " }
    15 = @{ D="15.28"; E="  -0.75%  " }
    16 = @{ D="2.664.35"; E="  +0.79%  " }
    17 = @{ D="2.319.77"; E="  +0.06%  " }
    18 = @{ D="42.712.45"; E="  +0.45%  " }
    19 = @{ D="7.45"; E="  -0.81%  " }
    20 = @{ E="  -0.42%  " }
    21 = @{ D="13.07"; E="  -12.35%  " }
    22 = @{ D="73.60"; E="  -0.49%  " }
    23 = @{ D="3.48"; E="  -1.70%  " }
    25 = @{ E="  +0.87%  " }
    26 = @{ D="1.00"; E="  -0.06%  " }
    27 = @{ E="  +11.95%  " }
    28 = @{ D="11.05"; E="  +1.05%  " }
    29 = @{ D="2.30"; E="  +0.88%  " }
    30 = @{ E="  +3.95%  " }
    31 = @{ D="22.43"; E="  -0.65%  " }
    32 = @{ D="166.21"; E="  +0.04%  " }
    33 = @{ D="0.0875"; E="  -0.45%  " }
    34 = @{ D="2.74"; E="  +5.04%  " }
    35 = @{ E="  -0.87%  " }
    36 = @{ D="4.71"; E="  +3.17%  " }
    37 = @{ E="  -1.92%  " }
    38 = @{ E="  +1.23%  " }
    39 = @{ D="2.82"; E="  +4.86%  " }
    40 = @{ D="3.68"; E="  -1.22%  " }
    41 = @{ D="1.60"; E="  +1.56%  " }
    42 = @{ D="104.28"; E="  +9.43%  " }
    43 = @{ D="70.71"; E="  +0.44%  " }
    44 = @{ E="  +1.69%  " }
    45 = @{ D="12.91"; E="  +4.81%  " }
    46 = @{ E="  -0.06%  " }
    47 = @{ D="112.56"; E="  -2.07%  " }
    48 = @{ B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.657.79"; E="  -2.91%  " }
    49 = @{ B="ordi"; C="https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"; D="76.97"; E="  -3.51%  " }
    50 = @{ E="  -0.15%  " }
    51 = @{ E="  +2.77%  " }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $cell = $ws.Cells.Item([int]$row, $colIndex[$col])
        # Force text storage so numeric-looking strings (e.g. "73.60",
        # "1.00") keep their original formatting/trailing zeros instead of
        # being coerced to a number.
        $cell.NumberFormat = "@"
        $cell.Value = $rowVals[$col]
        $cell.Style = "Normal"
    }
}
